# Classify instance files for better inspection:
# update the "AVG" (lower bound sample size) column values for several
# benchmark instances on sheet "04 Sep".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("04 Sep")

$ws.Range("E10").Value = 182
$ws.Range("E11").Value = 255
$ws.Range("E12").Value = 338
$ws.Range("E20").Value = 64
$ws.Range("E22").Value = 290

# Move the active selection, as recorded in the saved workbook.
$ws.Range("E23").Select()
